$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the data range first so numeric-looking strings
# (e.g. "596.54", "1.00") are stored as text, matching the source inlineStr cells,
# instead of being auto-converted to numbers by Excel.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "68.054.46"
$ws.Range("E2").Value = "  +1.00%  "
$ws.Range("D3").Value = "2.521.96"
$ws.Range("E3").Value = "  -1.20%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "596.54"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("D6").Value = "175.83"
$ws.Range("E6").Value = "  +0.83%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "0.529"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "2.519.73"
$ws.Range("E9").Value = "  -1.19%  "
$ws.Range("E10").Value = "  -0.48%  "
$ws.Range("E11").Value = "  +2.15%  "
$ws.Range("E12").Value = "  -1.34%  "
$ws.Range("E13").Value = "  -2.71%  "
$ws.Range("D14").Value = "26.57"
$ws.Range("E14").Value = "  -2.09%  "
$ws.Range("D15").Value = "2.980.06"
$ws.Range("E15").Value = "  -1.07%  "
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").Value = "67.815.70"
$ws.Range("E17").Value = "  +0.90%  "
$ws.Range("D18").Value = "2.552.21"
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("D19").Value = "12.01"
$ws.Range("E19").Value = "  +5.06%  "
$ws.Range("D20").Value = "8.06"
$ws.Range("E20").Value = "  -0.81%  "
$ws.Range("D21").Value = "364.76"
$ws.Range("E21").Value = "  +2.37%  "
$ws.Range("D22").Value = "4.17"
$ws.Range("E22").Value = "  -1.88%  "
$ws.Range("D23").Value = "4.65"
$ws.Range("E23").Value = "  -0.97%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "71.19"
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("E26").Value = "  -3.12%  "
$ws.Range("D27").Value = "10.12"
$ws.Range("E27").Value = "  +2.30%  "
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.15%  "
$ws.Range("D29").Value = "2.643.36"
$ws.Range("E29").Value = "  -1.07%  "
$ws.Range("D30").Value = "0.0₃0981"
$ws.Range("E30").Value = "  -2.13%  "
$ws.Range("D31").Value = "8.31"
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("D32").Value = "531.14"
$ws.Range("E32").Value = "  -1.35%  "
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").Value = "1.91"
$ws.Range("E33").Value = "  +2.28%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "1.33"
$ws.Range("E34").Value = "  -1.78%  "
$ws.Range("E35").Value = "  -2.58%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  -2.74%  "
$ws.Range("D38").Value = "156.53"
$ws.Range("E38").Value = "  -0.97%  "
$ws.Range("D39").Value = "18.74"
$ws.Range("E39").Value = "  -0.48%  "
$ws.Range("D40").Value = "18.68"
$ws.Range("E40").Value = "  +1.28%  "
$ws.Range("E41").Value = "  -0.53%  "
$ws.Range("D42").Value = "0.351"
$ws.Range("E42").Value = "  -2.04%  "
$ws.Range("D43").Value = "5.15"
$ws.Range("E43").Value = "  -1.24%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("E45").Value = "  -2.69%  "
$ws.Range("D46").Value = "147.37"
$ws.Range("E46").Value = "  -3.45%  "
$ws.Range("E47").Value = "  -1.86%  "
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("D49").Value = "0.0₆0275"
$ws.Range("E49").Value = "  -2.30%  "
$ws.Range("D50").Value = "1.72"
$ws.Range("E50").Value = "  -0.76%  "
$ws.Range("D51").Value = "0.0752"
$ws.Range("E51").Value = "  -1.29%  "

# Restore the default "Normal" style so cells keep the same (unstyled) look
# as the rest of the sheet instead of retaining the temporary text format.
$dataRange.Style = "Normal"
